$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue "D2" '60.587.14'
Set-TextValue "E2" '  +4.33%  '
Set-TextValue "D3" '2.337.08'
Set-TextValue "E3" '  +2.32%  '
Set-TextValue "E4" '  +0.02%  '
Set-TextValue "D5" '548.53'
Set-TextValue "E5" '  +2.66%  '
Set-TextValue "D6" '132.06'
Set-TextValue "E6" '  +0.74%  '
Set-TextValue "E7" '  -0.04%  '
Set-TextValue "E8" '  -0.79%  '
Set-TextValue "D9" '2.335.34'
Set-TextValue "E9" '  +2.27%  '
Set-TextValue "E10" '  +1.70%  '
Set-TextValue "D11" '5.52'
Set-TextValue "E11" '  +0.83%  '
Set-TextValue "E12" '  +0.16%  '
Set-TextValue "E13" '  +1.76%  '
Set-TextValue "D14" '23.89'
Set-TextValue "E14" '  +1.91%  '
Set-TextValue "D15" '2.754.15'
Set-TextValue "E15" '  +2.33%  '
Set-TextValue "D16" '60.522.29'
Set-TextValue "E16" '  +4.33%  '
Set-TextValue "D17" '0.0000134'
Set-TextValue "E17" '  +1.35%  '
Set-TextValue "D18" '2.332.67'
Set-TextValue "E18" '  +3.02%  '
Set-TextValue "D19" '10.66'
Set-TextValue "E19" '  +1.62%  '
Set-TextValue "E20" '  -0.21%  '
Set-TextValue "D21" '315.87'
Set-TextValue "E21" '  +0.84%  '
Set-TextValue "D22" '6.69'
Set-TextValue "E22" '  +4.12%  '
Set-TextValue "E23" '  -0.26%  '
Set-TextValue "D24" '64.29'
Set-TextValue "E24" '  +1.94%  '
Set-TextValue "E25" '  +1.43%  '
Set-TextValue "E26" '  +0.01%  '
Set-TextValue "D27" '7.88'
Set-TextValue "E27" '  -0.97%  '
Set-TextValue "D28" '1.36'
Set-TextValue "E28" '  +8.12%  '
Set-TextValue "E29" '  +13.24%  '
Set-TextValue "D30" '173.75'
Set-TextValue "E30" '  +1.83%  '
Set-TextValue "E31" '  +2.64%  '
Set-TextValue "D32" '0.0₃0740'
Set-TextValue "E32" '  +2.59%  '
Set-TextValue "E33" '  +3.71%  '
Set-TextValue "E34" '  +11.60%  '
Set-TextValue "D35" '0.382'
Set-TextValue "E35" '  +0.67%  '
Set-TextValue "D36" '17.97'
Set-TextValue "E36" '  +0.77%  '
Set-TextValue "E37" '  +0.02%  '
Set-TextValue "E38" '  +0.03%  '
Set-TextValue "E39" '  +5.57%  '
Set-TextValue "D40" '327.36'
Set-TextValue "E40" '  +14.40%  '
Set-TextValue "B41" 'Stacks'
Set-TextValue "C41" 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue "D41" '1.54'
Set-TextValue "E41" '  +3.32%  '
Set-TextValue "B42" 'OKB'
Set-TextValue "C42" 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue "D42" '38.16'
Set-TextValue "E42" '  -0.37%  '
Set-TextValue "D43" '140.08'
Set-TextValue "E43" '  +0.14%  '
Set-TextValue "E44" '  +1.66%  '
Set-TextValue "E45" '  -0.72%  '
Set-TextValue "D46" '19.43'
Set-TextValue "E46" '  +7.73%  '
Set-TextValue "D47" '0.0499'
Set-TextValue "E47" '  +1.14%  '
Set-TextValue "E48" '  +2.05%  '
Set-TextValue "E49" '  +21.37%  '
Set-TextValue "E50" '  +1.95%  '
Set-TextValue "D51" '11.03'
Set-TextValue "E51" '  +0.88%  '
